$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Förändrad" date for unchanged rows 2 and 3 ---
$ws.Cells.Item(2, 3).Value = 46065
$ws.Cells.Item(3, 3).Value = 46065

# Row 4: A 10617-2025
$ws.Cells.Item(4, 1).Value = 'A 10617-2025'
$ws.Cells.Item(4, 2).Value = 45721.61657407408
$ws.Cells.Item(4, 3).Value = 46065
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = 1.9
$ws.Cells.Item(4, 8).Value = 2
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 2
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = 'Lökgroda' + ([char]13+[char]10) + 'Strandpadda'
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 10617-2025 artfynd.xlsx", "A 10617-2025")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 10617-2025 karta.png", "A 10617-2025")'
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 10617-2025 FSC-klagomål.docx", "A 10617-2025")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 10617-2025 FSC-klagomål mail.docx", "A 10617-2025")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 10617-2025 tillsynsbegäran.docx", "A 10617-2025")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 10617-2025 tillsynsbegäran mail.docx", "A 10617-2025")'
$ws.Rows.Item(4).RowHeight = 15

# Row 5: A 26984-2022
$ws.Cells.Item(5, 1).Value = 'A 26984-2022'
$ws.Cells.Item(5, 2).Value = 44740.72620370371
$ws.Cells.Item(5, 3).Value = 46065
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = 2.9
$ws.Cells.Item(5, 8).Value = 2
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 2
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = 2
$ws.Cells.Item(5, 18).Value = 'Lökgroda' + ([char]13+[char]10) + 'Strandpadda'
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 26984-2022 artfynd.xlsx", "A 26984-2022")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 26984-2022 karta.png", "A 26984-2022")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 26984-2022 FSC-klagomål.docx", "A 26984-2022")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 26984-2022 FSC-klagomål mail.docx", "A 26984-2022")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 26984-2022 tillsynsbegäran.docx", "A 26984-2022")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 26984-2022 tillsynsbegäran mail.docx", "A 26984-2022")'
$ws.Rows.Item(5).RowHeight = 15

# Row 6: A 22776-2024
$ws.Cells.Item(6, 1).Value = 'A 22776-2024'
$ws.Cells.Item(6, 2).Value = 45448.47032407407
$ws.Cells.Item(6, 3).Value = 46065
$ws.Cells.Item(6, 6).Value = ""
$ws.Cells.Item(6, 7).Value = 5.6
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(6, 17).Value = 2
$ws.Cells.Item(6, 18).Value = 'Slåttergubbe' + ([char]13+[char]10) + 'Grönvit nattviol'
$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 22776-2024 artfynd.xlsx", "A 22776-2024")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 22776-2024 karta.png", "A 22776-2024")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 22776-2024 FSC-klagomål.docx", "A 22776-2024")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 22776-2024 FSC-klagomål mail.docx", "A 22776-2024")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 22776-2024 tillsynsbegäran.docx", "A 22776-2024")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 22776-2024 tillsynsbegäran mail.docx", "A 22776-2024")'
$ws.Rows.Item(6).RowHeight = 15

# Row 7: A 52355-2021
$ws.Cells.Item(7, 1).Value = 'A 52355-2021'
$ws.Cells.Item(7, 2).Value = 44464
$ws.Cells.Item(7, 3).Value = 46065
$ws.Cells.Item(7, 6).Value = ""
$ws.Cells.Item(7, 7).Value = 2.5
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 1
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 1
$ws.Cells.Item(7, 16).Value = 1
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = 'Skogsalm'
$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 52355-2021 artfynd.xlsx", "A 52355-2021")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 52355-2021 karta.png", "A 52355-2021")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 52355-2021 FSC-klagomål.docx", "A 52355-2021")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 52355-2021 FSC-klagomål mail.docx", "A 52355-2021")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 52355-2021 tillsynsbegäran.docx", "A 52355-2021")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 52355-2021 tillsynsbegäran mail.docx", "A 52355-2021")'
$ws.Rows.Item(7).RowHeight = 15

# Row 8: A 55145-2023
$ws.Cells.Item(8, 1).Value = 'A 55145-2023'
$ws.Cells.Item(8, 2).Value = 45237
$ws.Cells.Item(8, 3).Value = 46065
$ws.Cells.Item(8, 6).Value = 'Kommuner'
$ws.Cells.Item(8, 7).Value = 4.5
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = 'Grovticka'
$ws.Cells.Item(8, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 55145-2023 artfynd.xlsx", "A 55145-2023")'
$ws.Cells.Item(8, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 55145-2023 karta.png", "A 55145-2023")'
$ws.Cells.Item(8, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 55145-2023 FSC-klagomål.docx", "A 55145-2023")'
$ws.Cells.Item(8, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 55145-2023 FSC-klagomål mail.docx", "A 55145-2023")'
$ws.Cells.Item(8, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 55145-2023 tillsynsbegäran.docx", "A 55145-2023")'
$ws.Cells.Item(8, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 55145-2023 tillsynsbegäran mail.docx", "A 55145-2023")'
$ws.Rows.Item(8).RowHeight = 15

# Row 9: A 35443-2021
$ws.Cells.Item(9, 1).Value = 'A 35443-2021'
$ws.Cells.Item(9, 2).Value = 44385
$ws.Cells.Item(9, 3).Value = 46065
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(9, 7).Value = 4.2
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 1
$ws.Cells.Item(9, 16).Value = 1
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = 'Klubbfibbla'
$ws.Cells.Item(9, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 35443-2021 artfynd.xlsx", "A 35443-2021")'
$ws.Cells.Item(9, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 35443-2021 karta.png", "A 35443-2021")'
$ws.Cells.Item(9, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 35443-2021 FSC-klagomål.docx", "A 35443-2021")'
$ws.Cells.Item(9, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 35443-2021 FSC-klagomål mail.docx", "A 35443-2021")'
$ws.Cells.Item(9, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 35443-2021 tillsynsbegäran.docx", "A 35443-2021")'
$ws.Cells.Item(9, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 35443-2021 tillsynsbegäran mail.docx", "A 35443-2021")'
$ws.Rows.Item(9).RowHeight = 15

# Row 10: A 4746-2025
$ws.Cells.Item(10, 1).Value = 'A 4746-2025'
$ws.Cells.Item(10, 2).Value = 45688.46011574074
$ws.Cells.Item(10, 3).Value = 46065
$ws.Cells.Item(10, 6).Value = ""
$ws.Cells.Item(10, 7).Value = 1.1
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 1
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 'Svartvit flugsnappare'
$ws.Cells.Item(10, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 4746-2025 artfynd.xlsx", "A 4746-2025")'
$ws.Cells.Item(10, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 4746-2025 karta.png", "A 4746-2025")'
$ws.Cells.Item(10, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 4746-2025 FSC-klagomål.docx", "A 4746-2025")'
$ws.Cells.Item(10, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 4746-2025 FSC-klagomål mail.docx", "A 4746-2025")'
$ws.Cells.Item(10, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 4746-2025 tillsynsbegäran.docx", "A 4746-2025")'
$ws.Cells.Item(10, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 4746-2025 tillsynsbegäran mail.docx", "A 4746-2025")'
$ws.Cells.Item(10, 26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/fåglar/A 4746-2025 prioriterade fågelarter.docx", "A 4746-2025")'
$ws.Rows.Item(10).RowHeight = 15

# Row 11: A 4816-2025
$ws.Cells.Item(11, 1).Value = 'A 4816-2025'
$ws.Cells.Item(11, 2).Value = 45688.62052083333
$ws.Cells.Item(11, 3).Value = 46065
$ws.Cells.Item(11, 6).Value = ""
$ws.Cells.Item(11, 7).Value = 7.7
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = 'Strutbräken'
$ws.Cells.Item(11, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 4816-2025 artfynd.xlsx", "A 4816-2025")'
$ws.Cells.Item(11, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 4816-2025 karta.png", "A 4816-2025")'
$ws.Cells.Item(11, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 4816-2025 FSC-klagomål.docx", "A 4816-2025")'
$ws.Cells.Item(11, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 4816-2025 FSC-klagomål mail.docx", "A 4816-2025")'
$ws.Cells.Item(11, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 4816-2025 tillsynsbegäran.docx", "A 4816-2025")'
$ws.Cells.Item(11, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 4816-2025 tillsynsbegäran mail.docx", "A 4816-2025")'
$ws.Rows.Item(11).RowHeight = 15

# Row 12: A 31764-2023
$ws.Cells.Item(12, 1).Value = 'A 31764-2023'
$ws.Cells.Item(12, 2).Value = 45118
$ws.Cells.Item(12, 3).Value = 46065
$ws.Cells.Item(12, 6).Value = 'Kommuner'
$ws.Cells.Item(12, 7).Value = 1.2
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = 'Skorpdyna'
$ws.Cells.Item(12, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 31764-2023 artfynd.xlsx", "A 31764-2023")'
$ws.Cells.Item(12, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 31764-2023 karta.png", "A 31764-2023")'
$ws.Cells.Item(12, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 31764-2023 FSC-klagomål.docx", "A 31764-2023")'
$ws.Cells.Item(12, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 31764-2023 FSC-klagomål mail.docx", "A 31764-2023")'
$ws.Cells.Item(12, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 31764-2023 tillsynsbegäran.docx", "A 31764-2023")'
$ws.Cells.Item(12, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 31764-2023 tillsynsbegäran mail.docx", "A 31764-2023")'
$ws.Rows.Item(12).RowHeight = 15

# Row 13: A 11732-2024
$ws.Cells.Item(13, 1).Value = 'A 11732-2024'
$ws.Cells.Item(13, 2).Value = 45373.69222222222
$ws.Cells.Item(13, 3).Value = 46065
$ws.Cells.Item(13, 6).Value = ""
$ws.Cells.Item(13, 7).Value = 2.8
$ws.Cells.Item(13, 8).Value = 1
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 1
$ws.Cells.Item(13, 16).Value = 1
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = 'Dvärgjohannesört'
$ws.Cells.Item(13, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 11732-2024 artfynd.xlsx", "A 11732-2024")'
$ws.Cells.Item(13, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 11732-2024 karta.png", "A 11732-2024")'
$ws.Cells.Item(13, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 11732-2024 FSC-klagomål.docx", "A 11732-2024")'
$ws.Cells.Item(13, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 11732-2024 FSC-klagomål mail.docx", "A 11732-2024")'
$ws.Cells.Item(13, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 11732-2024 tillsynsbegäran.docx", "A 11732-2024")'
$ws.Cells.Item(13, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 11732-2024 tillsynsbegäran mail.docx", "A 11732-2024")'
$ws.Rows.Item(13).RowHeight = 15

# Row 14: A 72269-2021
$ws.Cells.Item(14, 1).Value = 'A 72269-2021'
$ws.Cells.Item(14, 2).Value = 44543
$ws.Cells.Item(14, 3).Value = 46065
$ws.Cells.Item(14, 6).Value = ""
$ws.Cells.Item(14, 7).Value = 3.2
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 1
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = 'Stor häxört'
$ws.Cells.Item(14, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 72269-2021 artfynd.xlsx", "A 72269-2021")'
$ws.Cells.Item(14, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 72269-2021 karta.png", "A 72269-2021")'
$ws.Cells.Item(14, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 72269-2021 FSC-klagomål.docx", "A 72269-2021")'
$ws.Cells.Item(14, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 72269-2021 FSC-klagomål mail.docx", "A 72269-2021")'
$ws.Cells.Item(14, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 72269-2021 tillsynsbegäran.docx", "A 72269-2021")'
$ws.Cells.Item(14, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 72269-2021 tillsynsbegäran mail.docx", "A 72269-2021")'
$ws.Rows.Item(14).RowHeight = 15

# Row 15: A 39260-2025
$ws.Cells.Item(15, 1).Value = 'A 39260-2025'
$ws.Cells.Item(15, 2).Value = 45889.37043981482
$ws.Cells.Item(15, 3).Value = 46065
$ws.Cells.Item(15, 6).Value = ""
$ws.Cells.Item(15, 7).Value = 7.6
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(15, 18).Value = 'Skogsödla'
$ws.Cells.Item(15, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 39260-2025 artfynd.xlsx", "A 39260-2025")'
$ws.Cells.Item(15, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 39260-2025 karta.png", "A 39260-2025")'
$ws.Cells.Item(15, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 39260-2025 FSC-klagomål.docx", "A 39260-2025")'
$ws.Cells.Item(15, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 39260-2025 FSC-klagomål mail.docx", "A 39260-2025")'
$ws.Cells.Item(15, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 39260-2025 tillsynsbegäran.docx", "A 39260-2025")'
$ws.Cells.Item(15, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 39260-2025 tillsynsbegäran mail.docx", "A 39260-2025")'
$ws.Rows.Item(15).RowHeight = 15

# Row 16: A 4792-2025
$ws.Cells.Item(16, 1).Value = 'A 4792-2025'
$ws.Cells.Item(16, 2).Value = 45688.57549768518
$ws.Cells.Item(16, 3).Value = 46065
$ws.Cells.Item(16, 6).Value = ""
$ws.Cells.Item(16, 7).Value = 5.9
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 1
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 15).Value = 1
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = 'Oxtungssvamp'
$ws.Cells.Item(16, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 4792-2025 artfynd.xlsx", "A 4792-2025")'
$ws.Cells.Item(16, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 4792-2025 karta.png", "A 4792-2025")'
$ws.Cells.Item(16, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 4792-2025 FSC-klagomål.docx", "A 4792-2025")'
$ws.Cells.Item(16, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 4792-2025 FSC-klagomål mail.docx", "A 4792-2025")'
$ws.Cells.Item(16, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 4792-2025 tillsynsbegäran.docx", "A 4792-2025")'
$ws.Cells.Item(16, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 4792-2025 tillsynsbegäran mail.docx", "A 4792-2025")'
$ws.Rows.Item(16).RowHeight = 15

# Row 17: A 59227-2025
$ws.Cells.Item(17, 1).Value = 'A 59227-2025'
$ws.Cells.Item(17, 2).Value = 45988.62253472222
$ws.Cells.Item(17, 3).Value = 46065
$ws.Cells.Item(17, 6).Value = ""
$ws.Cells.Item(17, 7).Value = 7.1
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = 'Scharlakansskål'
$ws.Cells.Item(17, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/artfynd/A 59227-2025 artfynd.xlsx", "A 59227-2025")'
$ws.Cells.Item(17, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/kartor/A 59227-2025 karta.png", "A 59227-2025")'
$ws.Cells.Item(17, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomål/A 59227-2025 FSC-klagomål.docx", "A 59227-2025")'
$ws.Cells.Item(17, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/klagomålsmail/A 59227-2025 FSC-klagomål mail.docx", "A 59227-2025")'
$ws.Cells.Item(17, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsyn/A 59227-2025 tillsynsbegäran.docx", "A 59227-2025")'
$ws.Cells.Item(17, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1281/tillsynsmail/A 59227-2025 tillsynsbegäran mail.docx", "A 59227-2025")'
$ws.Rows.Item(17).RowHeight = 15

# Row 18: A 22411-2021
$ws.Cells.Item(18, 1).Value = 'A 22411-2021'
$ws.Cells.Item(18, 2).Value = 44326
$ws.Cells.Item(18, 3).Value = 46065
$ws.Cells.Item(18, 6).Value = 'Kommuner'
$ws.Cells.Item(18, 7).Value = 3.4
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 0
$ws.Rows.Item(18).RowHeight = 15

# Row 19: A 68005-2021
$ws.Cells.Item(19, 1).Value = 'A 68005-2021'
$ws.Cells.Item(19, 2).Value = 44525
$ws.Cells.Item(19, 3).Value = 46065
$ws.Cells.Item(19, 6).Value = ""
$ws.Cells.Item(19, 7).Value = 0.5
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 0
$ws.Rows.Item(19).RowHeight = 15

# Row 20: A 56855-2021
$ws.Cells.Item(20, 1).Value = 'A 56855-2021'
$ws.Cells.Item(20, 2).Value = 44481
$ws.Cells.Item(20, 3).Value = 46065
$ws.Cells.Item(20, 6).Value = ""
$ws.Cells.Item(20, 7).Value = 2.7
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 0
$ws.Rows.Item(20).RowHeight = 15

# Row 21: A 21767-2023
$ws.Cells.Item(21, 1).Value = 'A 21767-2023'
$ws.Cells.Item(21, 2).Value = 45065
$ws.Cells.Item(21, 3).Value = 46065
$ws.Cells.Item(21, 6).Value = ""
$ws.Cells.Item(21, 7).Value = 1.9
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 0
$ws.Rows.Item(21).RowHeight = 15

# Row 22: A 45725-2022
$ws.Cells.Item(22, 1).Value = 'A 45725-2022'
$ws.Cells.Item(22, 2).Value = 44846
$ws.Cells.Item(22, 3).Value = 46065
$ws.Cells.Item(22, 6).Value = ""
$ws.Cells.Item(22, 7).Value = 1.2
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 0
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 0
$ws.Rows.Item(22).RowHeight = 15

# Row 23: A 32972-2024
$ws.Cells.Item(23, 1).Value = 'A 32972-2024'
$ws.Cells.Item(23, 2).Value = 45517.42064814815
$ws.Cells.Item(23, 3).Value = 46065
$ws.Cells.Item(23, 6).Value = 'Kommuner'
$ws.Cells.Item(23, 7).Value = 20.7
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(23, 15).Value = 0
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 0
$ws.Rows.Item(23).RowHeight = 15

# Row 24: A 33985-2021
$ws.Cells.Item(24, 1).Value = 'A 33985-2021'
$ws.Cells.Item(24, 2).Value = 44378
$ws.Cells.Item(24, 3).Value = 46065
$ws.Cells.Item(24, 6).Value = ""
$ws.Cells.Item(24, 7).Value = 1.5
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 0
$ws.Rows.Item(24).RowHeight = 15

# Row 25: A 5059-2025
$ws.Cells.Item(25, 1).Value = 'A 5059-2025'
$ws.Cells.Item(25, 2).Value = 45691.47211805556
$ws.Cells.Item(25, 3).Value = 46065
$ws.Cells.Item(25, 6).Value = ""
$ws.Cells.Item(25, 7).Value = 2.3
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 0
$ws.Rows.Item(25).RowHeight = 15

# Row 26: A 29075-2023
$ws.Cells.Item(26, 1).Value = 'A 29075-2023'
$ws.Cells.Item(26, 2).Value = 45104
$ws.Cells.Item(26, 3).Value = 46065
$ws.Cells.Item(26, 6).Value = ""
$ws.Cells.Item(26, 7).Value = 0.7
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = 0
$ws.Cells.Item(26, 14).Value = 0
$ws.Cells.Item(26, 15).Value = 0
$ws.Cells.Item(26, 16).Value = 0
$ws.Cells.Item(26, 17).Value = 0
$ws.Rows.Item(26).RowHeight = 15

# Row 27: A 20057-2025
$ws.Cells.Item(27, 1).Value = 'A 20057-2025'
$ws.Cells.Item(27, 2).Value = 45772
$ws.Cells.Item(27, 3).Value = 46065
$ws.Cells.Item(27, 6).Value = ""
$ws.Cells.Item(27, 7).Value = 1.3
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 0
$ws.Cells.Item(27, 14).Value = 0
$ws.Cells.Item(27, 15).Value = 0
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(27, 17).Value = 0
$ws.Rows.Item(27).RowHeight = 15

# Row 28: A 31773-2023
$ws.Cells.Item(28, 1).Value = 'A 31773-2023'
$ws.Cells.Item(28, 2).Value = 45118
$ws.Cells.Item(28, 3).Value = 46065
$ws.Cells.Item(28, 6).Value = 'Kommuner'
$ws.Cells.Item(28, 7).Value = 0.8
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(28, 14).Value = 0
$ws.Cells.Item(28, 15).Value = 0
$ws.Cells.Item(28, 16).Value = 0
$ws.Cells.Item(28, 17).Value = 0
$ws.Rows.Item(28).RowHeight = 15

# Row 29: A 13754-2022
$ws.Cells.Item(29, 1).Value = 'A 13754-2022'
$ws.Cells.Item(29, 2).Value = 44649
$ws.Cells.Item(29, 3).Value = 46065
$ws.Cells.Item(29, 6).Value = ""
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0
$ws.Rows.Item(29).RowHeight = 15

# Row 30: A 29247-2025
$ws.Cells.Item(30, 1).Value = 'A 29247-2025'
$ws.Cells.Item(30, 2).Value = 45824
$ws.Cells.Item(30, 3).Value = 46065
$ws.Cells.Item(30, 6).Value = ""
$ws.Cells.Item(30, 7).Value = 0.9
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(30, 14).Value = 0
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0
$ws.Rows.Item(30).RowHeight = 15

# Row 31: A 59918-2022
$ws.Cells.Item(31, 1).Value = 'A 59918-2022'
$ws.Cells.Item(31, 2).Value = 44909
$ws.Cells.Item(31, 3).Value = 46065
$ws.Cells.Item(31, 6).Value = ""
$ws.Cells.Item(31, 7).Value = 1.6
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(31, 14).Value = 0
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(31, 17).Value = 0
$ws.Rows.Item(31).RowHeight = 15

# Row 32: A 44646-2023
$ws.Cells.Item(32, 1).Value = 'A 44646-2023'
$ws.Cells.Item(32, 2).Value = 45189
$ws.Cells.Item(32, 3).Value = 46065
$ws.Cells.Item(32, 6).Value = ""
$ws.Cells.Item(32, 7).Value = 3.2
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = 0
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(32, 17).Value = 0
$ws.Rows.Item(32).RowHeight = 15

# Row 33: A 15984-2025
$ws.Cells.Item(33, 1).Value = 'A 15984-2025'
$ws.Cells.Item(33, 2).Value = 45749
$ws.Cells.Item(33, 3).Value = 46065
$ws.Cells.Item(33, 6).Value = ""
$ws.Cells.Item(33, 7).Value = 2.4
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = 0
$ws.Cells.Item(33, 14).Value = 0
$ws.Cells.Item(33, 15).Value = 0
$ws.Cells.Item(33, 16).Value = 0
$ws.Cells.Item(33, 17).Value = 0
$ws.Rows.Item(33).RowHeight = 15

# Row 34: A 29083-2023
$ws.Cells.Item(34, 1).Value = 'A 29083-2023'
$ws.Cells.Item(34, 2).Value = 45105
$ws.Cells.Item(34, 3).Value = 46065
$ws.Cells.Item(34, 6).Value = ""
$ws.Cells.Item(34, 7).Value = 0.5
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(34, 14).Value = 0
$ws.Cells.Item(34, 15).Value = 0
$ws.Cells.Item(34, 16).Value = 0
$ws.Cells.Item(34, 17).Value = 0
$ws.Rows.Item(34).RowHeight = 15

# Row 35: A 49297-2023
$ws.Cells.Item(35, 1).Value = 'A 49297-2023'
$ws.Cells.Item(35, 2).Value = 45210
$ws.Cells.Item(35, 3).Value = 46065
$ws.Cells.Item(35, 6).Value = ""
$ws.Cells.Item(35, 7).Value = 1
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = 0
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(35, 17).Value = 0
$ws.Rows.Item(35).RowHeight = 15

# Row 36: A 34484-2025
$ws.Cells.Item(36, 1).Value = 'A 34484-2025'
$ws.Cells.Item(36, 2).Value = 45847.39517361111
$ws.Cells.Item(36, 3).Value = 46065
$ws.Cells.Item(36, 6).Value = 'Kommuner'
$ws.Cells.Item(36, 7).Value = 1.2
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(36, 17).Value = 0
$ws.Rows.Item(36).RowHeight = 15

# Row 37: A 14058-2025
$ws.Cells.Item(37, 1).Value = 'A 14058-2025'
$ws.Cells.Item(37, 2).Value = 45740.40390046296
$ws.Cells.Item(37, 3).Value = 46065
$ws.Cells.Item(37, 6).Value = ""
$ws.Cells.Item(37, 7).Value = 2.7
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = 0
$ws.Cells.Item(37, 14).Value = 0
$ws.Cells.Item(37, 15).Value = 0
$ws.Cells.Item(37, 16).Value = 0
$ws.Cells.Item(37, 17).Value = 0
$ws.Rows.Item(37).RowHeight = 15

# Row 38: A 28885-2022
$ws.Cells.Item(38, 1).Value = 'A 28885-2022'
$ws.Cells.Item(38, 2).Value = 44749
$ws.Cells.Item(38, 3).Value = 46065
$ws.Cells.Item(38, 6).Value = ""
$ws.Cells.Item(38, 7).Value = 2.4
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 0
$ws.Cells.Item(38, 14).Value = 0
$ws.Cells.Item(38, 15).Value = 0
$ws.Cells.Item(38, 16).Value = 0
$ws.Cells.Item(38, 17).Value = 0
$ws.Rows.Item(38).RowHeight = 15

# Row 39: A 16499-2025
$ws.Cells.Item(39, 1).Value = 'A 16499-2025'
$ws.Cells.Item(39, 2).Value = 45751.60246527778
$ws.Cells.Item(39, 3).Value = 46065
$ws.Cells.Item(39, 6).Value = ""
$ws.Cells.Item(39, 7).Value = 0.5
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = 0
$ws.Cells.Item(39, 14).Value = 0
$ws.Cells.Item(39, 15).Value = 0
$ws.Cells.Item(39, 16).Value = 0
$ws.Cells.Item(39, 17).Value = 0
$ws.Rows.Item(39).RowHeight = 15

# Row 40: A 59733-2023
$ws.Cells.Item(40, 1).Value = 'A 59733-2023'
$ws.Cells.Item(40, 2).Value = 45257
$ws.Cells.Item(40, 3).Value = 46065
$ws.Cells.Item(40, 6).Value = ""
$ws.Cells.Item(40, 7).Value = 8.9
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = 0
$ws.Cells.Item(40, 15).Value = 0
$ws.Cells.Item(40, 16).Value = 0
$ws.Cells.Item(40, 17).Value = 0
$ws.Rows.Item(40).RowHeight = 15

# Row 41: A 55158-2023
$ws.Cells.Item(41, 1).Value = 'A 55158-2023'
$ws.Cells.Item(41, 2).Value = 45237
$ws.Cells.Item(41, 3).Value = 46065
$ws.Cells.Item(41, 6).Value = 'Kommuner'
$ws.Cells.Item(41, 7).Value = 0.7
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = 0
$ws.Cells.Item(41, 14).Value = 0
$ws.Cells.Item(41, 15).Value = 0
$ws.Cells.Item(41, 16).Value = 0
$ws.Cells.Item(41, 17).Value = 0
$ws.Rows.Item(41).RowHeight = 15

# Row 42: A 8601-2023
$ws.Cells.Item(42, 1).Value = 'A 8601-2023'
$ws.Cells.Item(42, 2).Value = 44977
$ws.Cells.Item(42, 3).Value = 46065
$ws.Cells.Item(42, 6).Value = ""
$ws.Cells.Item(42, 7).Value = 1.5
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(42, 17).Value = 0
$ws.Rows.Item(42).RowHeight = 15

# Row 43: A 17519-2023
$ws.Cells.Item(43, 1).Value = 'A 17519-2023'
$ws.Cells.Item(43, 2).Value = 45036
$ws.Cells.Item(43, 3).Value = 46065
$ws.Cells.Item(43, 6).Value = ""
$ws.Cells.Item(43, 7).Value = 0.5
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(43, 14).Value = 0
$ws.Cells.Item(43, 15).Value = 0
$ws.Cells.Item(43, 16).Value = 0
$ws.Cells.Item(43, 17).Value = 0
$ws.Rows.Item(43).RowHeight = 15

# Row 44: A 6034-2025
$ws.Cells.Item(44, 1).Value = 'A 6034-2025'
$ws.Cells.Item(44, 2).Value = 45695.64231481482
$ws.Cells.Item(44, 3).Value = 46065
$ws.Cells.Item(44, 6).Value = ""
$ws.Cells.Item(44, 7).Value = 10.2
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(44, 17).Value = 0
$ws.Rows.Item(44).RowHeight = 15

# Row 45: A 2248-2025
$ws.Cells.Item(45, 1).Value = 'A 2248-2025'
$ws.Cells.Item(45, 2).Value = 45673.49231481482
$ws.Cells.Item(45, 3).Value = 46065
$ws.Cells.Item(45, 6).Value = ""
$ws.Cells.Item(45, 7).Value = 11.7
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(45, 17).Value = 0
$ws.Rows.Item(45).RowHeight = 15

# Row 46: A 58337-2025
$ws.Cells.Item(46, 1).Value = 'A 58337-2025'
$ws.Cells.Item(46, 2).Value = 45985.48895833334
$ws.Cells.Item(46, 3).Value = 46065
$ws.Cells.Item(46, 6).Value = ""
$ws.Cells.Item(46, 7).Value = 2
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 0
$ws.Cells.Item(46, 14).Value = 0
$ws.Cells.Item(46, 15).Value = 0
$ws.Cells.Item(46, 16).Value = 0
$ws.Cells.Item(46, 17).Value = 0
$ws.Rows.Item(46).RowHeight = 15

# Row 47: A 58488-2025
$ws.Cells.Item(47, 1).Value = 'A 58488-2025'
$ws.Cells.Item(47, 2).Value = 45985.66210648148
$ws.Cells.Item(47, 3).Value = 46065
$ws.Cells.Item(47, 6).Value = ""
$ws.Cells.Item(47, 7).Value = 9.5
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0
$ws.Cells.Item(47, 14).Value = 0
$ws.Cells.Item(47, 15).Value = 0
$ws.Cells.Item(47, 16).Value = 0
$ws.Cells.Item(47, 17).Value = 0
$ws.Rows.Item(47).RowHeight = 15

# Row 48: A 58335-2025
$ws.Cells.Item(48, 1).Value = 'A 58335-2025'
$ws.Cells.Item(48, 2).Value = 45985.48332175926
$ws.Cells.Item(48, 3).Value = 46065
$ws.Cells.Item(48, 6).Value = ""
$ws.Cells.Item(48, 7).Value = 2.5
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 0
$ws.Cells.Item(48, 14).Value = 0
$ws.Cells.Item(48, 15).Value = 0
$ws.Cells.Item(48, 16).Value = 0
$ws.Cells.Item(48, 17).Value = 0
$ws.Rows.Item(48).RowHeight = 15

# Row 49: A 59206-2025
$ws.Cells.Item(49, 1).Value = 'A 59206-2025'
$ws.Cells.Item(49, 2).Value = 45988.60881944445
$ws.Cells.Item(49, 3).Value = 46065
$ws.Cells.Item(49, 6).Value = ""
$ws.Cells.Item(49, 7).Value = 6.9
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = 0
$ws.Cells.Item(49, 14).Value = 0
$ws.Cells.Item(49, 15).Value = 0
$ws.Cells.Item(49, 16).Value = 0
$ws.Cells.Item(49, 17).Value = 0
$ws.Rows.Item(49).RowHeight = 15

